$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- sheet1 (rri) view: deselect tab, change selection/scroll ---
$ws1.Range("B14:D32").Select()

# --- add the new "pduration" sheet after "rri" ---
$newws = $wb.Worksheets.Add([Type]::Missing, $ws1)
$newws.Name = "pduration"

# --- header row ---
$newws.Range("A1").Value = "rate"
$newws.Range("B1").Value = "pv"
$newws.Range("C1").Value = "fv"
$newws.Range("D1").Value = "pduration"

# --- data rows (A:C values, D formula added after) ---
$newws.Range("A2").Value = 0.015309470499731193
$newws.Range("B2").Value = -5
$newws.Range("C2").Value = -6
$newws.Range("A3").Value = -1
$newws.Range("B3").Value = -5
$newws.Range("C3").Value = 0
$newws.Range("A4").Value = 0
$newws.Range("B4").Value = -1
$newws.Range("C4").Value = -1
$newws.Range("A5").Value = 0
$newws.Range("B5").Value = 300
$newws.Range("C5").Value = 300
$newws.Range("A6").Value = 0.1
$newws.Range("B6").Value = 0
$newws.Range("C6").Value = 100
$newws.Range("A7").Value = 0.1
$newws.Range("B7").Value = 100
$newws.Range("C7").Value = 0
$newws.Range("A8").Value = 0.02426318074098921
$newws.Range("B8").Value = 300
$newws.Range("C8").Value = 400
$newws.Range("A9").Value = 0.24092317318260137
$newws.Range("B9").Value = 300
$newws.Range("C9").Value = 4000
$newws.Range("A10").Value = 0.50341274654387536
$newws.Range("B10").Value = 300
$newws.Range("C10").Value = 40000
$newws.Range("A11").Value = 0.012058882052318642
$newws.Range("B11").Value = 300
$newws.Range("C11").Value = 400
$newws.Range("A12").Value = 0.11396731243901459
$newws.Range("B12").Value = 300
$newws.Range("C12").Value = 4000
$newws.Range("A13").Value = 0.22613732776711237
$newws.Range("B13").Value = 300
$newws.Range("C13").Value = 40000
$newws.Range("A14").Value = 0.0075993101546305564
$newws.Range("B14").Value = 300
$newws.Range("C14").Value = 400
$newws.Range("A15").Value = 0.070541853470322824
$newws.Range("B15").Value = 300
$newws.Range("C15").Value = 4000
$newws.Range("A16").Value = 0.13741628093790048
$newws.Range("B16").Value = 300
$newws.Range("C16").Value = 40000
$newws.Range("A17").Value = 0.98822504304098735
$newws.Range("B17").Value = 10000
$newws.Range("C17").Value = 2441880
$newws.Range("A18").Value = 0.046635139392105618
$newws.Range("B18").Value = 5000
$newws.Range("C18").Value = 6000
$newws.Range("A19").Value = 0.18920711500272103
$newws.Range("B19").Value = 5000
$newws.Range("C19").Value = 10000
$newws.Range("A20").Value = 0.10000000000000009
$newws.Range("B20").Value = 250
$newws.Range("C20").Value = 275
$newws.Range("A21").Value = 0.41421356237309492
$newws.Range("B21").Value = 250
$newws.Range("C21").Value = 500
$newws.Range("A22").Value = 0.5211809843045565
$newws.Range("B22").Value = 250
$newws.Range("C22").Value = 880
$newws.Range("A23").Value = 0.025000000000000001
$newws.Range("B23").Value = 2000
$newws.Range("C23").Value = 2200
$newws.Range("A24").Formula = "=0.025/12"
$newws.Range("B24").Value = 1000
$newws.Range("C24").Value = 1200

# --- D column formulas: D2 stand-alone, D3:D24 as one shared-fill formula ---
$newws.Range("D2").Formula = "=_xlfn.PDURATION(A2,B2,C2)"
$newws.Range("D3:D24").Formula = "=_xlfn.PDURATION(A3,B3,C3)"

# --- format column D width like the original ---
$newws.Columns.Item(4).ColumnWidth = 11.08984375

# --- create Table2 over the new sheet data ---
$lo = $newws.ListObjects.Add(1, $newws.Range("A1:D24"), [Type]::Missing, 1)
$lo.Name = "Table2"
$lo.TableStyle = "TableStyleMedium7"

# --- select/activate final view state on the new sheet ---
$newws.Activate()
$newws.Range("A24").Select()
